# Build the ATM test-case table on Sheet1.
# Reproduces: header row, 8 test-case rows (A:test#, B:description,
# C:initial balance, D:user input, E:expected output), plus formatting
# (bold numbers, centered text, currency format, wrapped/centered long
# text, bold+italic header) and column widths / row heights / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. Values — written in the same order the author would naturally have
#    typed them (header, then description column, then user-input
#    column, then the last row, then the expected-output column) so
#    that the workbook's shared-string table comes out in that order.
# ---------------------------------------------------------------------

# Header row
$ws.Range("A1").Value = "Test Case Number"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Initial Balance"
$ws.Range("D1").Value = "User Input"
$ws.Range("E1").Value = "Expected Output"

# Test case # column (A) and Initial Balance column (C) for all 8 rows
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 3).Value = 1000
}

# Description column (B) for the first 7 rows
$descriptions = @(
    "Valid Deposit",
    "Valid Withdrawal",
    "Withdraw more than balance",
    "Negative Deposit",
    "Negative Withdrawal",
    "View Balance",
    "Invalid Input"
)
for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $descriptions[$i]
}

# User Input column (D) for the first 7 rows
$userInputs = @("D, 200", "W, 100", "W, 1100", "D, -100", "W, -50", "V", "X")
for ($i = 0; $i -lt $userInputs.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $userInputs[$i]
}

# Last row (row 9 = test case 8): description then user input
$ws.Range("B9").Value = "Exit Program"
$ws.Range("D9").Value = "E"

# Expected Output column (E) for all 8 rows
$expectedOutputs = @(
    "Balance: `$1200.00",
    "Balance: `$900.00",
    "Error: Insufficient funds.",
    "Error: Cannot deposit a negative amount.",
    "Error: Cannot withdraw a negative amount.",
    "Balance: `$1000.00",
    "Error: Invalid choice.",
    '"Thank you for using the ATM. Goodbye!"'
)
for ($i = 0; $i -lt $expectedOutputs.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $expectedOutputs[$i]
}

# ---------------------------------------------------------------------
# 2. Formatting — each distinct cell style is built once on a single
#    cell (so the style is created cleanly) and then copied as
#    format-only onto the rest of the cells that share it.
# ---------------------------------------------------------------------

# Style 1: column E baseline — wrap text + vertical centered (no
# horizontal centering). Applied to the bottom helper cell E10.
$ws.Range("E10").WrapText = $true
$ws.Range("E10").VerticalAlignment = $xlCenter

# Style 2: wrap text + fully centered — the long error/expected-output
# messages in E4:E9.
$base2 = $ws.Range("E4")
$base2.WrapText = $true
$base2.HorizontalAlignment = $xlCenter
$base2.VerticalAlignment = $xlCenter
$base2.Copy()
$ws.Range("E5:E9").PasteSpecial($xlPasteFormats)

# Style 3: centered, no wrap — Description (B) and User Input (D)
# columns, plus the short expected-output cells E2:E3.
$base3 = $ws.Range("B2")
$base3.HorizontalAlignment = $xlCenter
$base3.VerticalAlignment = $xlCenter
$base3.Copy()
$dest3 = $ws.Range("B3:B9,D2:D9,E2:E3")
$dest3.PasteSpecial($xlPasteFormats)

# Style 4: currency number format + centered — Initial Balance (C).
$base4 = $ws.Range("C2")
$base4.NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'
$base4.HorizontalAlignment = $xlCenter
$base4.VerticalAlignment = $xlCenter
$base4.Copy()
$ws.Range("C3:C9").PasteSpecial($xlPasteFormats)

# Style 5: bold + centered — Test Case Number column (A).
$base5 = $ws.Range("A2")
$base5.Font.Bold = $true
$base5.HorizontalAlignment = $xlCenter
$base5.VerticalAlignment = $xlCenter
$base5.Copy()
$ws.Range("A3:A9").PasteSpecial($xlPasteFormats)

# Style 6: bold + italic + centered — header row.
$base6 = $ws.Range("A1")
$base6.Font.Bold = $true
$base6.Font.Italic = $true
$base6.HorizontalAlignment = $xlCenter
$base6.VerticalAlignment = $xlCenter
$base6.Copy()
$ws.Range("B1:E1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Column widths / row heights
# ---------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 18.125
$ws.Columns.Item(2).ColumnWidth = 24.75
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 37.25

$rowHeights = @(18.75, 30, 29.25, 27.75, 25.5, 24.75, 23.25, 28.5, 27.75)
for ($i = 0; $i -lt $rowHeights.Length; $i++) {
    $ws.Rows.Item($i + 1).RowHeight = $rowHeights[$i]
}

# ---------------------------------------------------------------------
# 4. Selection — matches the author's cursor position when they saved.
# ---------------------------------------------------------------------

$ws.Range("H8").Select() | Out-Null
